$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 120; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
